$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hourly crypto price/volume refresh. Price cells in column D are stored as
# literal text (e.g. "67.806.36", "1.00") so that values which otherwise look
# numeric do not get silently reinterpreted as numbers by Excel; the cell is
# briefly switched to Text format to force that, then restored to the default
# "Normal" style so no visible formatting change is left behind.

$ws.Range("D2").Value = "67.806.36"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.297.35"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "186.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "581.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.129"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.407"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "3.873.44"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.135"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "67.934.85"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000167"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "3.294.45"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "447.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +12.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.16%  "
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "3.471.23"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.514"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000118"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.188"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.784"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("D43").Value = "2.692.53"
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0673"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "324.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0276"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "31.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.985"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.60%  "
